$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content of the used range so the shared-strings table gets
# rebuilt fresh, in the exact order cells are (re)written below.
$ws.Range("A1:P3").ClearContents()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "Width-1 (mm)"
$ws.Range("C1").Value = "Width-2 (mm)"
$ws.Range("D1").Value = "Extra Width Generated (mm)"
$ws.Range("E1").Value = "Trim (mm)"
$ws.Range("F1").Value = "Mother Coil Width (mm)"
$ws.Range("G1").Value = "Grade"
$ws.Range("H1").Value = "Mother Coil Weight (kg)"
$ws.Range("I1").Value = "Mother Coil Length (m)"
$ws.Range("J1").Value = "Weight-1 (kg)"
$ws.Range("K1").Value = "Weight-2 (kg)"
$ws.Range("L1").Value = "Extra Weight (kg)"
$ws.Range("M1").Value = "Width-3 (mm)"
$ws.Range("N1").Value = "Width-4 (mm)"
$ws.Range("O1").Value = "Weight-3 (kg)"
$ws.Range("P1").Value = "Weight-4 (kg)"

# --- Row 2 (data row, id = 0) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 180
$ws.Range("C2").Value = 190
$ws.Range("D2").Value = 320
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 700
$ws.Range("G2").Value = "M75_BAO"
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 822.6491771451605
$ws.Range("J2").Value = 257.1428571428571
$ws.Range("K2").Value = 271.4285714285714
$ws.Range("L2").Value = 457.1428571428572
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0

# --- Row 3 (data row, id = 1) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 110
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 140
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 710
$ws.Range("G3").Value = "M75_BAO"
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1216.593853524533
$ws.Range("J3").Value = 232.394366197183
$ws.Range("K3").Value = 274.6478873239437
$ws.Range("L3").Value = 295.7746478873239
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 170
$ws.Range("O3").Value = 316.9014084507042
$ws.Range("P3").Value = 359.1549295774648
